# Update stats for 2025-09 (row 22 in Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B22").Value = 6297
$ws.Range("D22").Value = 5855817
$ws.Range("E22").Value = 929.9375893282515
$ws.Range("F22").Value = 8.400757445343432
$ws.Range("H22").Value = 27.34577110100123
